{"js": "// Apply the nanny-info text updates via Office.js (Word JavaScript API).\nconst body = context.document.body;\n\nconst replacements = [\n  { find: \"This morning, Joey woke up at 07:00 AM.\", replace: \"This morning, Joey woke up at 10:34 AM.\" },\n  { find: \"Her first nap should be at 09:00 AM.\", replace: \"Her first nap should be at 10:34 AM.\" },\n  { find: \"For lunch today, we have Cherries.\", replace: \"For lunch today, we have Dhbdvdbd.\" },\n  { find: \"For dinner today, we have Berries.\", replace: \"For dinner today, we have Sbn end.\" },\n  { find: \"As a reminder Don\\u2019t mess up.\", replace: \"As a reminder Dbnejdbcf.\" }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the nanny-info text updates via Word COM interop.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"This morning, Joey woke up at 07:00 AM.\"; Replace = \"This morning, Joey woke up at 10:34 AM.\" },\n    @{ Find = \"Her first nap should be at 09:00 AM.\"; Replace = \"Her first nap should be at 10:34 AM.\" },\n    @{ Find = \"For lunch today, we have Cherries.\"; Replace = \"For lunch today, we have Dhbdvdbd.\" },\n    @{ Find = \"For dinner today, we have Berries.\"; Replace = \"For dinner today, we have Sbn end.\" },\n    @{ Find = \"As a reminder Don\u2019t mess up.\"; Replace = \"As a reminder Dbnejdbcf.\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n"}
